# Chiffres COVID-19 Valais — daily data update
#
# The "Cumul cas positifs" (col B), "Total hospitalisations" (col H) and
# "Cumul deces" (col J/K) columns are TODAY()-driven running totals, so we
# only need to type the new raw daily figures into the input columns
# (C, E, F, G, L, M) for the days that were actually reported; Excel's
# own recalculation (which `run_com` performs automatically once this
# script finishes) fills in the dependent B/H/J/K cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- correction to an already-reported day -------------------------------
# "Nb nouveaux cas positifs" for 2021-05-17 (row 506) was corrected
# 29 -> 28; every subsequent cumulative total in column B cascades by -1
# automatically through the shared formula.
$ws.Range("C506").Value2 = 28

# --- two more corrected/completed days -----------------------------------
$ws.Range("C600").Value2 = 24
$ws.Range("C601").Value2 = 61

# --- three newly reported days (2021-10-19, 10-20, 10-21) ----------------
$ws.Range("C602").Value2 = 25
$ws.Range("E602").Value2 = 2
$ws.Range("F602").Value2 = 1
$ws.Range("G602").Value2 = 9
$ws.Range("L602").Value2 = 0
$ws.Range("M602").Value2 = 0

$ws.Range("C603").Value2 = 52
$ws.Range("E603").Value2 = 2
$ws.Range("F603").Value2 = 1
$ws.Range("G603").Value2 = 8
$ws.Range("L603").Value2 = 0
$ws.Range("M603").Value2 = 0

$ws.Range("C604").Value2 = 1
$ws.Range("E604").Value2 = 2
$ws.Range("F604").Value2 = 1
$ws.Range("G604").Value2 = 8
$ws.Range("L604").Value2 = 0
$ws.Range("M604").Value2 = 0
